$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(175, 5).Value = 60
$ws.Cells.Item(175, 11).Value = 60.61

$ws.Cells.Item(202, 5).Value = 30
$ws.Cells.Item(202, 11).Value = 30.3

$ws.Cells.Item(483, 5).Value = 15
$ws.Cells.Item(483, 11).Value = 55.56

$ws.Cells.Item(510, 5).Value = 2
$ws.Cells.Item(510, 11).Value = 7.41

$ws.Cells.Item(795, 5).Value = 52
$ws.Cells.Item(795, 11).Value = 44.44

$ws.Cells.Item(822, 5).Value = 8
$ws.Cells.Item(822, 11).Value = 6.84

$ws.Cells.Item(180, 5).Value = 2
$ws.Cells.Item(180, 7).Value = 99
$ws.Cells.Item(180, 8).Value = 0
$ws.Cells.Item(180, 9).Value = 0
$ws.Cells.Item(180, 10).Value = 99
$ws.Cells.Item(180, 11).Value = 2.02

$ws.Cells.Item(182, 5).Value = 2
$ws.Cells.Item(182, 7).Value = 99
$ws.Cells.Item(182, 8).Value = 0
$ws.Cells.Item(182, 9).Value = 0
$ws.Cells.Item(182, 10).Value = 99
$ws.Cells.Item(182, 11).Value = 2.02

$ws.Cells.Item(183, 5).Value = 35
$ws.Cells.Item(183, 7).Value = 99
$ws.Cells.Item(183, 8).Value = 0
$ws.Cells.Item(183, 9).Value = 0
$ws.Cells.Item(183, 10).Value = 99
$ws.Cells.Item(183, 11).Value = 35.35

$ws.Cells.Item(207, 5).Value = 4
$ws.Cells.Item(207, 7).Value = 99
$ws.Cells.Item(207, 8).Value = 0
$ws.Cells.Item(207, 9).Value = 0
$ws.Cells.Item(207, 10).Value = 99
$ws.Cells.Item(207, 11).Value = 4.04

$ws.Cells.Item(209, 5).Value = 7
$ws.Cells.Item(209, 7).Value = 99
$ws.Cells.Item(209, 8).Value = 0
$ws.Cells.Item(209, 9).Value = 0
$ws.Cells.Item(209, 10).Value = 99
$ws.Cells.Item(209, 11).Value = 7.07

$ws.Cells.Item(210, 5).Value = 58
$ws.Cells.Item(210, 7).Value = 99
$ws.Cells.Item(210, 8).Value = 0
$ws.Cells.Item(210, 9).Value = 0
$ws.Cells.Item(210, 10).Value = 99
$ws.Cells.Item(210, 11).Value = 58.59

$ws.Cells.Item(490, 5).Value = 1
$ws.Cells.Item(490, 7).Value = 27
$ws.Cells.Item(490, 8).Value = 0
$ws.Cells.Item(490, 9).Value = 0
$ws.Cells.Item(490, 10).Value = 27
$ws.Cells.Item(490, 11).Value = 3.7

$ws.Cells.Item(491, 5).Value = 11
$ws.Cells.Item(491, 7).Value = 27
$ws.Cells.Item(491, 8).Value = 0
$ws.Cells.Item(491, 9).Value = 0
$ws.Cells.Item(491, 10).Value = 27
$ws.Cells.Item(491, 11).Value = 40.74

$ws.Cells.Item(517, 5).Value = 2
$ws.Cells.Item(517, 7).Value = 27
$ws.Cells.Item(517, 8).Value = 0
$ws.Cells.Item(517, 9).Value = 0
$ws.Cells.Item(517, 10).Value = 27
$ws.Cells.Item(517, 11).Value = 7.41

$ws.Cells.Item(518, 5).Value = 23
$ws.Cells.Item(518, 7).Value = 27
$ws.Cells.Item(518, 8).Value = 0
$ws.Cells.Item(518, 9).Value = 0
$ws.Cells.Item(518, 10).Value = 27
$ws.Cells.Item(518, 11).Value = 85.19

$ws.Cells.Item(802, 5).Value = 7
$ws.Cells.Item(802, 7).Value = 117
$ws.Cells.Item(802, 8).Value = 0
$ws.Cells.Item(802, 9).Value = 0
$ws.Cells.Item(802, 10).Value = 117
$ws.Cells.Item(802, 11).Value = 5.98

$ws.Cells.Item(803, 5).Value = 58
$ws.Cells.Item(803, 7).Value = 117
$ws.Cells.Item(803, 8).Value = 0
$ws.Cells.Item(803, 9).Value = 0
$ws.Cells.Item(803, 10).Value = 117
$ws.Cells.Item(803, 11).Value = 49.57

$ws.Cells.Item(827, 5).Value = 1
$ws.Cells.Item(827, 7).Value = 117
$ws.Cells.Item(827, 8).Value = 0
$ws.Cells.Item(827, 9).Value = 0
$ws.Cells.Item(827, 10).Value = 117
$ws.Cells.Item(827, 11).Value = 0.85

$ws.Cells.Item(829, 5).Value = 13
$ws.Cells.Item(829, 7).Value = 117
$ws.Cells.Item(829, 8).Value = 0
$ws.Cells.Item(829, 9).Value = 0
$ws.Cells.Item(829, 10).Value = 117
$ws.Cells.Item(829, 11).Value = 11.11

$ws.Cells.Item(830, 5).Value = 95
$ws.Cells.Item(830, 7).Value = 117
$ws.Cells.Item(830, 8).Value = 0
$ws.Cells.Item(830, 9).Value = 0
$ws.Cells.Item(830, 10).Value = 117
$ws.Cells.Item(830, 11).Value = 81.2

for ($r = 620; $r -le 935; $r++) {
    $ws.Cells.Item($r, 1).Value = "Wadjemup"
}
